$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A4').Value = 'CS402[AP]  /  []'
$ws.Range('B4').Value = 'B401[PD]  /  []'
$ws.Range('C4').Value = 'M401[SMa]  /  []'
$ws.Range('D4').Value = 'CS401[AH]  /  []'
$ws.Range('E4').Value = 'Free Period!'
$ws.Range('F4').Value = 'CS403[MDu]  /  []'
$ws.Range('G4').Value = 'CH401[SC]  /  []'
$ws.Range('A6').Value = 'CS491[AH, MDu]  /  CS492[AP, BDu]'
$ws.Range('B6').Value = 'CS491[AH, MDu]  /  CS492[AP, BDu]'
$ws.Range('C6').Value = 'CS491[AH, MDu]  /  CS492[AP, BDu]'
$ws.Range('D6').Value = 'CS491[AH, MDu]  /  CS492[AP, BDu]'
$ws.Range('E6').Value = 'B401[PD]  /  []'
$ws.Range('F6').Value = 'Free Period!'
$ws.Range('G6').Value = 'CS402[AP]  /  []'
$ws.Range('A8').Value = 'CS402[AP]  /  []'
$ws.Range('B8').Value = 'M401[SMa]  /  []'
$ws.Range('C8').Value = 'CS401[AH]  /  []'
$ws.Range('D8').Value = 'CS492[AP, BDu]  /  CS491[AH, MDu]'
$ws.Range('E8').Value = 'CS492[AP, BDu]  /  CS491[AH, MDu]'
$ws.Range('F8').Value = 'CS492[AP, BDu]  /  CS491[AH, MDu]'
$ws.Range('G8').Value = 'CS492[AP, BDu]  /  CS491[AH, MDu]'
$ws.Range('A10').Value = 'Free Period!'
$ws.Range('B10').Value = 'B401[PD]  /  []'
$ws.Range('C10').Value = 'CH401[SC]  /  []'
$ws.Range('D10').Value = 'Free Period!'
$ws.Range('E10').Value = 'Free Period!'
$ws.Range('F10').Value = 'CS403[MDu]  /  []'
$ws.Range('G10').Value = 'Free Period!'
$ws.Range('A12').Value = 'CS401[AH]  /  []'
$ws.Range('B12').Value = 'M401(T)[DC, GY]  /  []'
$ws.Range('C12').Value = 'CH401[SC]  /  []'
$ws.Range('D12').Value = 'M401[SMa]  /  []'
$ws.Range('E12').Value = 'Free Period!'
$ws.Range('F12').Value = 'CS403[MDu]  /  []'
$ws.Range('G12').Value = 'Free Period!'
$ws.Range('A16').Value = 'Free Period!'
$ws.Range('B16').Value = 'IT403[SU]  /  []'
$ws.Range('C16').Value = 'IT401[RCh]  /  []'
$ws.Range('D16').Value = 'Free Period!'
$ws.Range('E16').Value = 'CH401[SC]  /  []'
$ws.Range('F16').Value = 'M401[ARC]  /  []'
$ws.Range('G16').Value = 'B401[PD]  /  []'
$ws.Range('A18').Value = 'B401[PD]  /  []'
$ws.Range('B18').Value = 'Free Period!'
$ws.Range('C18').Value = 'IT402[SSR]  /  []'
$ws.Range('D18').Value = 'IT492[RCh, RG]  /  IT491[SSR, KDa]'
$ws.Range('E18').Value = 'IT492[RCh, RG]  /  IT491[SSR, KDa]'
$ws.Range('F18').Value = 'IT492[RCh, RG]  /  IT491[SSR, KDa]'
$ws.Range('G18').Value = 'IT492[RCh, RG]  /  IT491[SSR, KDa]'
$ws.Range('A20').Value = 'IT402[SSR]  /  []'
$ws.Range('B20').Value = 'IT403[SU]  /  []'
$ws.Range('C20').Value = 'IT401[RCh]  /  []'
$ws.Range('D20').Value = 'Free Period!'
$ws.Range('E20').Value = 'CH401[SC]  /  []'
$ws.Range('F20').Value = 'M401[ARC]  /  []'
$ws.Range('G20').Value = 'Free Period!'
$ws.Range('A22').Value = 'CH401[SC]  /  []'
$ws.Range('B22').Value = 'IT403[SU]  /  []'
$ws.Range('C22').Value = 'IT401[RCh]  /  []'
$ws.Range('D22').Value = 'Free Period!'
$ws.Range('E22').Value = 'Free Period!'
$ws.Range('F22').Value = 'Free Period!'
$ws.Range('G22').Value = 'Free Period!'
$ws.Range('A24').Value = 'M401[ARC]  /  []'
$ws.Range('B24').Value = 'B401[PD]  /  []'
$ws.Range('C24').Value = 'IT402[SSR]  /  []'
$ws.Range('D24').Value = 'IT491[SSR, KDa]  /  IT492[RCh, RG]'
$ws.Range('E24').Value = 'IT491[SSR, KDa]  /  IT492[RCh, RG]'
$ws.Range('F24').Value = 'IT491[SSR, KDa]  /  IT492[RCh, RG]'
$ws.Range('G24').Value = 'IT491[SSR, KDa]  /  IT492[RCh, RG]'
$ws.Range('A28').Value = 'Free Period!'
$ws.Range('B28').Value = 'ECE404[SDe]  /  []'
$ws.Range('C28').Value = 'M401[SLa]  /  []'
$ws.Range('D28').Value = 'Free Period!'
$ws.Range('E28').Value = 'B401[PD]  /  []'
$ws.Range('F28').Value = 'ECE401[SMC]  /  []'
$ws.Range('G28').Value = 'ECE403[PP]  /  []'
$ws.Range('A30').Value = 'ECE402[SG]  /  []'
$ws.Range('B30').Value = 'ECE491[PC, SMC]  /  ECE492[DK, JA]'
$ws.Range('C30').Value = 'ECE491[PC, SMC]  /  ECE492[DK, JA]'
$ws.Range('D30').Value = 'ECE491[PC, SMC]  /  ECE492[DK, JA]'
$ws.Range('E30').Value = 'HU491[KB]  /  M491[SLa, SRC]'
$ws.Range('F30').Value = 'HU491[KB]  /  M491[SLa, SRC]'
$ws.Range('G30').Value = 'ECE404[SDe]  /  []'
$ws.Range('A32').Value = 'ECE493[PP, BC]  /  ECE491[PC, SMC]'
$ws.Range('B32').Value = 'ECE493[PP, BC]  /  ECE491[PC, SMC]'
$ws.Range('C32').Value = 'ECE493[PP, BC]  /  ECE491[PC, SMC]'
$ws.Range('D32').Value = 'ECE401[SMC]  /  []'
$ws.Range('E32').Value = 'ECE492[DK, JA]  /  ECE493[PP, BC]'
$ws.Range('F32').Value = 'ECE492[DK, JA]  /  ECE493[PP, BC]'
$ws.Range('G32').Value = 'ECE492[DK, JA]  /  ECE493[PP, BC]'
$ws.Range('A34').Value = 'M491[SLa, SRC]  /  HU491[KB]'
$ws.Range('B34').Value = 'M491[SLa, SRC]  /  HU491[KB]'
$ws.Range('C34').Value = 'ECE404[SDe]  /  []'
$ws.Range('D34').Value = 'ECE402[SG]  /  []'
$ws.Range('E34').Value = 'B401[PD]  /  []'
$ws.Range('F34').Value = 'M401[SLa]  /  []'
$ws.Range('G34').Value = 'ECE403[PP]  /  []'
$ws.Range('A36').Value = 'Free Period!'
$ws.Range('B36').Value = 'Free Period!'
$ws.Range('C36').Value = 'M401[SLa]  /  []'
$ws.Range('D36').Value = 'ECE402[SG]  /  []'
$ws.Range('E36').Value = 'B401[PD]  /  []'
$ws.Range('F36').Value = 'ECE401[SMC]  /  []'
$ws.Range('G36').Value = 'ECE403[PP]  /  []'
$ws.Range('A40').Value = 'CH401[PD]  /  []'
$ws.Range('B40').Value = 'EE492[RND, JA]  /  EE494[BDC, NCS]'
$ws.Range('C40').Value = 'EE492[RND, JA]  /  EE494[BDC, NCS]'
$ws.Range('D40').Value = 'EE492[RND, JA]  /  EE494[BDC, NCS]'
$ws.Range('E40').Value = 'EE401[BDC]  /  []'
$ws.Range('F40').Value = 'EE403[KR]  /  []'
$ws.Range('G40').Value = 'EE404[ArD]  /  []'
$ws.Range('A42').Value = 'EE401[BDC]  /  []'
$ws.Range('B42').Value = 'HU401[ACh]  /  []'
$ws.Range('C42').Value = 'EE403[KR]  /  []'
$ws.Range('D42').Value = 'EE404[ArD]  /  []'
$ws.Range('E42').Value = 'Free Period!'
$ws.Range('F42').Value = 'EE402[RND]  /  []'
$ws.Range('G42').Value = 'Free Period!'
$ws.Range('A44').Value = 'EE401[BDC]  /  []'
$ws.Range('B44').Value = 'EE494[BDC, NCS]  /  EE493[ArD, SDG]'
$ws.Range('C44').Value = 'EE494[BDC, NCS]  /  EE493[ArD, SDG]'
$ws.Range('D44').Value = 'EE494[BDC, NCS]  /  EE493[ArD, SDG]'
$ws.Range('E44').Value = 'Free Period!'
$ws.Range('F44').Value = 'EE402[RND]  /  []'
$ws.Range('G44').Value = 'HU401[ACh]  /  []'
$ws.Range('A46').Value = 'Free Period!'
$ws.Range('B46').Value = 'EE491[KR, IB]  /  EE492[RND, JA]'
$ws.Range('C46').Value = 'EE491[KR, IB]  /  EE492[RND, JA]'
$ws.Range('D46').Value = 'EE491[KR, IB]  /  EE492[RND, JA]'
$ws.Range('E46').Value = 'Free Period!'
$ws.Range('F46').Value = 'CH401[PD]  /  []'
$ws.Range('G46').Value = 'HU401[ACh]  /  []'
$ws.Range('A48').Value = 'EE403[KR]  /  []'
$ws.Range('B48').Value = 'EE493[ArD, SDG]  /  EE491[KR, IB]'
$ws.Range('C48').Value = 'EE493[ArD, SDG]  /  EE491[KR, IB]'
$ws.Range('D48').Value = 'EE493[ArD, SDG]  /  EE491[KR, IB]'
$ws.Range('E48').Value = 'EE402[RND]  /  []'
$ws.Range('F48').Value = 'CH401[PD]  /  []'
$ws.Range('G48').Value = 'EE404[ArD]  /  []'
